$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S1_Generation")
for ($r=9; $r -le 9; $r++) {
    for ($c=1; $c -le 8; $c++) {
        $v = $ws.Cells.Item($r, $c).Value()
        Write-Host "$r,$c : $v"
    }
}
for ($r=16; $r -le 16; $r++) {
    for ($c=1; $c -le 8; $c++) {
        $v = $ws.Cells.Item($r, $c).Value()
        Write-Host "$r,$c : $v"
    }
}
